# Add a new worksheet "2021_non_res" at the end of the workbook containing
# the 2019-2021 non-resident abortion counts by state, matching the
# "add non-resident data" commit.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet after the last existing sheet -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "2021_non_res"

# --- Populate data -----------------------------------------------------
# Fill values in the same order the original author typed them so that new
# shared-string entries land in the expected order (abortions, Kentucky,
# Illinois, Michigan, Tennessee, state - with "Ohio"/"Other" reusing
# already-existing shared strings).
$newSheet.Range("B1").Value = "abortions"

$newSheet.Range("A2").Value = "Kentucky"
$newSheet.Range("B2").Value = 264

$newSheet.Range("A3").Value = "Illinois"
$newSheet.Range("B3").Value = 71

$newSheet.Range("A4").Value = "Michigan"
$newSheet.Range("B4").Value = 56

$newSheet.Range("A5").Value = "Ohio"
$newSheet.Range("B5").Value = 40

$newSheet.Range("A6").Value = "Tennessee"
$newSheet.Range("B6").Value = 11

$newSheet.Range("A7").Value = "Other"
$newSheet.Range("B7").Value = 23

$newSheet.Range("A1").Value = "state"

# --- Selections on other touched sheets ---------------------------------
# 2021_gestation_weeks selection moves from C1 to F13 and loses tabSelected.
$gw = $wb.Worksheets.Item("2021_gestation_weeks")
$gw.Range("F13").Select()

# --- Make the new sheet the active / selected tab last -------------------
# so it ends up tabSelected in the xml and workbookView.activeTab points at
# it (2021_monthly, which previously held tabSelected, loses it
# automatically once another sheet becomes active).
$newSheet.Activate()
$newSheet.Range("E4").Select()
